# ---------------------------------------------------------------------------
# "Add files via upload" - reshape the single-sheet run-log workbook into
# three sheets: accept (trimmed original), reject (bigger log incl. a
# blank-row gap), and font_line (a third, smaller log).
# ---------------------------------------------------------------------------

# Helper: write a value into a cell while defending against Excel's
# automatic text->date / text->number coercion for values that merely look
# like a date ("06-05-2023") or a long digit run ("8858998585076"). Any
# value that starts with a digit gets a leading quote-prefix, which forces
# Excel to store it as literal text without the quote becoming part of the
# stored value.
function Set-TextCell($ws, $row, $col, $val) {
    if ($val -match '^[0-9]') {
        $ws.Cells.Item($row, $col).Value = "'" + $val
    } else {
        $ws.Cells.Item($row, $col).Value = $val
    }
}

function Set-DataRow($ws, $row, $a, $b, $c) {
    Set-TextCell $ws $row 1 $a
    Set-TextCell $ws $row 2 $b
    Set-TextCell $ws $row 3 $c
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the existing (only) sheet to "accept" and trim it down to the
#    header + first three data rows (dimension A1:C4).
# ---------------------------------------------------------------------------
$accept = $wb.Worksheets.Item(1)
$accept.Name = "accept"
$accept.Rows("5:11").Delete()

# ---------------------------------------------------------------------------
# 2. Insert "reject" right after "accept", and "font_line" right after
#    "reject" - keeps final tab order accept, reject, font_line.
# ---------------------------------------------------------------------------
$reject = $wb.Worksheets.Add($null, $accept)
$reject.Name = "reject"

$fontLine = $wb.Worksheets.Add($null, $reject)
$fontLine.Name = "font_line"

# ---------------------------------------------------------------------------
# 3. Populate "reject" - header + 9 rows, a blank row 11/12 gap, then one
#    more row at row 13.
# ---------------------------------------------------------------------------
Set-DataRow $reject 1 "d/m/y" "time" "run"
Set-DataRow $reject 2 "06-05-2023" "20:41:08" "th0000-0"
Set-DataRow $reject 3 "06-05-2023" "20:41:08" "th0000-1"
Set-DataRow $reject 4 "06-05-2023" "20:41:09" "th0000-2"
Set-DataRow $reject 5 "06-05-2023" "20:41:09" "th0000-3"
Set-DataRow $reject 6 "06-05-2023" "20:41:09" "th0000-4"
Set-DataRow $reject 7 "06-05-2024" "20:41:10" "th0000-5"
Set-DataRow $reject 8 "06-05-2025" "20:41:11" "th0000-6"
Set-DataRow $reject 9 "06-05-2026" "20:41:12" "th0000-7"
Set-DataRow $reject 10 "06-05-2023" "20:41:08" "th0000-8"
Set-DataRow $reject 13 "25-05-2023" "16:07:25" "8858998585076"

$reject.Columns.Item(1).ColumnWidth = 10.71
$reject.Columns.Item(2).ColumnWidth = 12.86
$reject.Columns.Item(3).ColumnWidth = 26.29
$reject.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 4. Populate "font_line" - header + two matching rows + one final row.
# ---------------------------------------------------------------------------
Set-DataRow $fontLine 1 "d/m/y" "time" "run"
Set-DataRow $fontLine 2 "24-05-2023" "13:25:04" "th0000-0"
Set-DataRow $fontLine 3 "24-05-2023" "13:25:04" "th0000-1"
Set-DataRow $fontLine 4 "25-05-2023" "16:02:43" "8858998585076"

$fontLine.Columns.Item(1).ColumnWidth = 13.29
$fontLine.Columns.Item(2).ColumnWidth = 18.29
$fontLine.Columns.Item(3).ColumnWidth = 22.14
$fontLine.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 5. View state: "accept" zoomed to 115% with A5 pre-selected, "font_line"
#    with A4:C4 selected, and "reject" left as the final active/selected
#    tab with A10 selected (matches activeTab=1 / tabSelected on sheet2).
# ---------------------------------------------------------------------------
$accept.Activate()
$excel.ActiveWindow.Zoom = 115
$accept.Range("A5").Select()

$fontLine.Activate()
$fontLine.Range("A4:C4").Select()

$reject.Activate()
$reject.Range("A10").Select()
